$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("H3").Value = 0.8247232472324724
$ws.Range("I3").Value = 0.01711175020542317
$ws.Range("J3").Value = 0.9
$ws.Range("K3").Value = 88.8

$ws.Range("Q3").Value = 33
$ws.Range("R3").Value = 61
$ws.Range("S3").Value = 86
$ws.Range("T3").Value = 92
$ws.Range("U3").Value = 115

$ws.Range("V3").Value = 4835
$ws.Range("W3").Value = 4807
$ws.Range("X3").Value = 4782
$ws.Range("Y3").Value = 4776
$ws.Range("Z3").Value = 4753

$ws.Range("AF3").Value = 0.993221
$ws.Range("AG3").Value = 0.987469
$ws.Range("AH3").Value = 0.982334
$ws.Range("AI3").Value = 0.981101
$ws.Range("AJ3").Value = 0.976376
